$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-8: "v" -> "VV" in columns F and G
$rows4to8 = 4,5,6,7,8
foreach ($r in $rows4to8) {
    $ws.Range("F$r").Value = "VV"
    $ws.Range("G$r").Value = "VV"
}

# Rows 22-24: add new F/G cells with "VV"
$rows22to24 = 22,23,24
foreach ($r in $rows22to24) {
    $ws.Range("F$r").Value = "VV"
    $ws.Range("G$r").Value = "VV"
}

# Row 25: add new H cell with "PDFJS-DIST"
$ws.Range("H25").Value = "PDFJS-DIST"

# Row 26: add new F/G cells with "VV"
$ws.Range("F26").Value = "VV"
$ws.Range("G26").Value = "VV"

# Update selection to match the new view (E21)
$ws.Range("E21").Select()
